# "Small words adjustment in Presentation"
#
# Slide 3 ("Why?") holds a SmartArt diagram (shape "Diagram 3") whose leaf
# nodes include the captions "I'm a rebel" and "I'm depressed". Rename them
# to "I'm happy" and "I'm sad" via the SmartArt object model so both the
# diagram data part (dataN.xml) and the cached diagram drawing part
# (drawingN.xml) pick up the new wording.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$diagramShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasSmartArt) {
        $diagramShape = $candidate
    }
}

$nodes = $diagramShape.SmartArt.AllNodes

for ($i = 1; $i -le $nodes.Count; $i++) {
    $node = $nodes.Item($i)
    $nodeText = $node.TextFrame.TextRange.Text
    if ($nodeText -eq "I’m a rebel") {
        $node.TextFrame.TextRange.Text = "I’m happy"
    }
    if ($nodeText -eq "I’m depressed") {
        $node.TextFrame.TextRange.Text = "I’m sad"
    }
}
